$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    0.093714853270462584,
    -0.0059999999140387672,
    -0.0039999999311568502,
    -0.0079999998724549215,
    -0.0029999999402230415,
    -0.0019999999482021025,
    -0.0099999998350721597,
    -0.0099999998352560127,
    -0.0019999999509172639,
    -0.0019999999543767188,
    -0.0029999999403882427,
    -0.0034999999352622879,
    0.0042319541705655439,
    -0.0079999998831894459,
    0.03411309792077688,
    -0.0019999999703741445,
    -0.0019999999686026726,
    -0.0039999999400555097,
    -0.0039999999419708665,
    -0.003999999940768717,
    -0.003999999941082244,
    0.0095951682939476868,
    -0.0049999999115737381,
    -0.01999999968994981,
    -0.019999999685413883,
    -0.0024999999408592544,
    -0.002499999940439146,
    -0.0019999999452551265,
    -0.0069999998737273472,
    -0.059999999126776338,
    0.01647415548313802,
    0.050146151347254531,
    -0.0039999999357274163
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $values[$i]
}
